# Generate Report for Archive
# Update the localization status report: files 867d690f-79e0-4e19-8a41-8176ded1bde6
# and 8b8a1244-c258-4e38-b606-24645a2cbf2a move from "Ready for handoff" to
# "In Translation" across the Overview summary sheet and the per-locale
# (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) for rows 3 & 4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"

# --- zh-cn sheet: Status column (C) for rows 3 & 4 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C4").Value = "In Translation"

# --- de-de sheet: Status column (C) for rows 3 & 4 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C4").Value = "In Translation"
